$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("Globo",  "Bom Dia Inter",       "Limpeza Pública", "2025-04-02T13:48", "Negativo", "TESTETESTETESTETESTETESTETESTETESTETESTE"),
    @("Record", "Balanço Geral",       "Limpeza Pública", "2025-04-02T13:48", "Neutro",   "TESTETESTETESTETESTETESTETESTETESTETESTE"),
    @("Record", "RJ No Ar TV Record",  "Codemca",         "2025-04-02T13:48", "Negativo", "TESTETESTETESTETESTETESTETESTETESTETESTETESTETESTETESTETESTE")
)

$startRow = 37
for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]
    for ($c = 1; $c -le $data.Length; $c++) {
        $ws.Cells.Item($r, $c).Value = $data[$c - 1]
    }
}
